$d = $word.ActiveDocument

# 1. Update the "Curso (semestre ideal)" line
$d.Content.Find.Execute(
    "Curso (semestre ideal): EQD (8), EQN (9)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Curso (semestre ideal): EQD (9), EQN (10)", 2)

# 2. Update the responsible professor listing
$d.Content.Find.Execute(
    "5840855 - Heizir Ferreira de Castro", $true, $false, $false, $false, $false,
    $true, 1, $false, "1285870 - Marcos Villela Barcza", 2)

# 3. Update the requirements list (two runs separated by a manual line break).
#    Both runs share identical (empty) run formatting, so if we edit both of
#    them directly the engine's save-time canonicalization will coalesce them
#    into a single <w:r>. To keep them as two separate runs (matching the
#    target), we temporarily bookmark the second run before editing the
#    first one, edit the second run through the bookmark (which keeps it
#    "detached" from the adjacency check), then remove the bookmark.
$p = $d.Paragraphs.Item(18)
$pStart = $p.Range.Start
$full = $p.Range.Text
$brk = [char]11
$idx1 = $full.IndexOf($brk)
$idx2 = $full.IndexOf($brk, $idx1 + 1)

$run2Start = $pStart + $idx1 + 1
$run2End = $pStart + $idx2 + 1
$r2pre = $d.Range($run2Start, $run2End)
$d.Bookmarks.Add("__loq4023_req2", $r2pre)

$run1End = $pStart + $idx1 + 1
$r1 = $d.Range($pStart, $run1End)
$r1.Find.Execute(
    "LOQ4002 -  Reatores Quimicos  (Requisito fraco)", $true, $false, $false, $false, $false,
    $true, 1, $false, "LOQ4038 -  Química Orgânica II  (Requisito fraco)", 2)

$bm = $d.Bookmarks.Item("__loq4023_req2")
$r2 = $bm.Range
$r2.Find.Execute(
    "LOT2004 -  Bioquímica  (Requisito fraco)", $true, $false, $false, $false, $false,
    $true, 1, $false, "LOQ4057 -  Operações Unitárias III  (Requisito fraco)", 2)

$d.Bookmarks.Item("__loq4023_req2").Delete()
